$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B42 currently stores "3" as text (a pre-existing data quirk). The new row 43
# needs that same "numeric-looking text" quirk, so copy B42's value+format down
# to B43 first (this preserves the plain, unstyled text representation),
# before overwriting B42 with a proper numeric 3.
$ws.Cells.Item(42, 2).Copy()
$ws.Cells.Item(43, 2).PasteSpecial(-4163)

# Fix row 42, column B: should now be the numeric value 3 (not text).
$ws.Cells.Item(42, 2).Value = 3

# Append the new annotation row 43.
$ws.Cells.Item(43, 1).Value = "Ruilin"
$ws.Cells.Item(43, 3).Value = "无"
$ws.Cells.Item(43, 4).Value = "QSN"
$ws.Cells.Item(43, 5).Value = "MET"
$ws.Cells.Item(43, 6).Value = "77ff87fb-cfc5-44ac-a4b7-cb33b05fed6f"
$ws.Cells.Item(43, 7).Value = "ByQpn1ZA-_annotated.xlsx"
$ws.Cells.Item(43, 8).Value = "If we know the regularization is fundamentally and mathematically wrong, why do we investigate its performance?"
